$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

$ws.Range("B2").Value = 10.37627124123017
$ws.Range("B3").Value = 0.12
$ws.Range("B4").Value = 1850.937149452345
$ws.Range("B5").Value = 21738.69177262505
$ws.Range("B8").Value = 1.056396007190265
$ws.Range("B9").Value = 0.6809103347581049
$ws.Range("B10").Value = 2.734361650964236
$ws.Range("B11").Value = 0.3427635586474372
